# Updated cryptos list (price + 1h volume refresh), matching the
# "Updated cryptos list on Wed Apr 26 03:57:07 UTC 2023 with GitHub Actions"
# commit. D-column price cells that look like plain numbers are written
# with a leading apostrophe so Excel stores them as text (preserving
# formats such as trailing zeros / multi-dot grouped numbers) instead of
# silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.364.35'
$ws.Range("E2").Value = '  +3.28%  '
$ws.Range("D3").Value = '1.866.31'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '''338.89'
$ws.Range("E5").Value = '  +1.90%  '
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = '''0.4681'
$ws.Range("E7").Value = '  +1.62%  '
$ws.Range("D8").Value = '''0.3959'
$ws.Range("E8").Value = '  +4.30%  '
$ws.Range("D9").Value = '''47.28'
$ws.Range("E9").Value = '  +1.66%  '
$ws.Range("D10").Value = '''0.08000'
$ws.Range("D11").Value = '''0.9989'
$ws.Range("E11").Value = '  +2.86%  '
$ws.Range("D12").Value = '''21.85'
$ws.Range("E12").Value = '  +4.27%  '
$ws.Range("E13").Value = '  +1.95%  '
$ws.Range("D14").Value = '''5.999'
$ws.Range("E14").Value = '  +1.93%  '
$ws.Range("D15").Value = '''7.229'
$ws.Range("E15").Value = '  +3.16%  '
$ws.Range("D16").Value = '''91.12'
$ws.Range("E16").Value = '  +3.82%  '
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").Value = '''0.00001040'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").Value = '''0.06620'
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").Value = '''17.49'
$ws.Range("E20").Value = '  +3.31%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = '28.375.18'
$ws.Range("E22").Value = '  +3.38%  '
$ws.Range("D23").Value = '''5.456'
$ws.Range("E23").Value = '  +2.39%  '
$ws.Range("E24").Value = '  +2.21%  '
$ws.Range("E25").Value = '  -1.37%  '
$ws.Range("D26").Value = '2.091.43'
$ws.Range("E26").Value = '  +1.61%  '
$ws.Range("E27").Value = '  +2.05%  '
$ws.Range("D28").Value = '''19.73'
$ws.Range("E28").Value = '  +2.22%  '
$ws.Range("D29").Value = '''2.121'
$ws.Range("E29").Value = '  +2.83%  '
$ws.Range("D30").Value = '''5.486'
$ws.Range("E30").Value = '  +3.36%  '
$ws.Range("D31").Value = '''120.22'
$ws.Range("D32").Value = '''0.9683'
$ws.Range("E32").Value = '  +1.83%  '
$ws.Range("E33").Value = '  +2.01%  '
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").Value = '''5.343'
$ws.Range("E35").Value = '  +2.44%  '
$ws.Range("D36").Value = '''1.369'
$ws.Range("E36").Value = '  +3.86%  '
$ws.Range("D37").Value = '''0.06085'
$ws.Range("E37").Value = '  +2.69%  '
$ws.Range("D38").Value = '''0.02239'
$ws.Range("E38").Value = '  +2.62%  '
$ws.Range("D39").Value = '''8.381'
$ws.Range("E39").Value = '  +3.97%  '
$ws.Range("D40").Value = '''1.184'
$ws.Range("E40").Value = '  +2.16%  '
$ws.Range("D41").Value = '''0.5936'
$ws.Range("E41").Value = '  +2.48%  '
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("E43").Value = '  +1.98%  '
$ws.Range("D44").Value = '''10.35'
$ws.Range("E44").Value = '  +3.50%  '
$ws.Range("E45").Value = '  +4.00%  '
$ws.Range("D46").Value = '''0.5566'
$ws.Range("E46").Value = '  +1.65%  '
$ws.Range("D47").Value = '''12.11'
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("D48").Value = '''1.953'
$ws.Range("E48").Value = '  +4.94%  '
$ws.Range("D49").Value = '''0.06855'
$ws.Range("E49").Value = '  +2.99%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '''111.39'
$ws.Range("E50").Value = '  +1.62%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '''2.044'
$ws.Range("E51").Value = '  +16.20%  '
